$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 0.1825419310453436
$ws.Range("C4").ClearContents()
$ws.Range("E4").Value = -0.1800933741311850
$ws.Range("C6").Value = 0.348613976222456
$ws.Range("C8").Value = -0.1384957661262676
$ws.Range("C9").Value = 1.692932643509826
$ws.Range("C10").Value = 1.566479473280191
$ws.Range("C11").Value = 1.020829760720687
$ws.Range("C12").Value = 0.7307568962937161
$ws.Range("E13").Value = 1.258913537332873
$ws.Range("C14").Value = 0.8188188121642126
$ws.Range("E14").Value = 0.960760217268164
$ws.Range("C15").Value = 1.019715257608933
$ws.Range("E16").Value = 1.375398114243231
$ws.Range("C17").Value = 2.173959184500363
$ws.Range("E17").Value = 1.566646323486065
$ws.Range("C18").Value = 1.9846842782967
$ws.Range("E18").Value = 1.47327408793585
$ws.Range("C19").Value = 1.707434489470039
$ws.Range("E20").Value = 1.681032827388385
$ws.Range("C21").Value = 1.456988786619817
$ws.Range("E21").Value = 1.84279714442821
$ws.Range("E23").Value = 1.657737120813474
$ws.Range("E24").Value = 1.580042106786372
$ws.Range("C25").Value = 1.260396653238427
$ws.Range("C26").Value = 1.064321453542272
$ws.Range("E26").Value = 0.7767182380207682
$ws.Range("E28").Value = 1.604795846351492
$ws.Range("E29").Value = 1.242807488305697
$ws.Range("C30").Value = 1.361817904277718
$ws.Range("E31").Value = 0.9049590709689692
$ws.Range("E32").Value = -1.215549235925817
$ws.Range("C33").Value = -7.03958082960261
$ws.Range("E33").Value = -9.851708704716611
$ws.Range("C34").Value = -4.352425014431327
$ws.Range("E34").Value = 0.9348518890383906
$ws.Range("C35").Value = -2.824222064391535
$ws.Range("E35").Value = -1.61400258701867
$ws.Range("C36").Value = -5.665308402785508
$ws.Range("E36").Value = -4.458023117238186
$ws.Range("E37").Value = -1.507094401446352
$ws.Range("E38").Value = 5.161235657134755
$ws.Range("C39").Value = 2.828271820504513
$ws.Range("E40").Value = -0.5534294478199198
$ws.Range("E41").Value = 0.7337317298176549
$ws.Range("E42").Value = 2.430255857698516
$ws.Range("C43").Value = 1.970682684899994
$ws.Range("E43").Value = 5.462250257438317
$ws.Range("C44").Value = -1.17492083522599
$ws.Range("E44").Value = 0.1140263184959744
$ws.Range("C46").Value = -0.9008525709169657
$ws.Range("E46").Value = 1.982587461121321
$ws.Range("C47").Value = 1.058598619486273
$ws.Range("C48").Value = 0.02017133142706573
$ws.Range("E48").Value = -0.3342090768663986
$ws.Range("C49").Value = 0.5138342970629317
$ws.Range("C50").Value = 0.2738544794132602
$ws.Range("C51").Value = 0.2382702494847733
$ws.Range("E51").Value = -0.01234339085524061
$ws.Range("E52").Value = -0.1189552196680155
